{"js": "// Replace the 25 division-problem answers in the table, in document order,\n// with their new values from the commit. The table cell paragraphs (plus the\n// leading date paragraph) are enumerated via context.document.body.paragraphs,\n// which walks the body in document order including paragraphs nested in\n// table cells (index 0 is the \"2024-02-20 Tuesday\" heading, indices 1-25 are\n// the populated answer cells interleaved with the blank spacer-row cells).\n\nconst replacements = [\n  \"34\u00f77=4, 6\",\n  \"86\u00f79=9, 5\",\n  \"90\u00f74=22, 2\",\n  \"85\u00f76=14, 1\",\n  \"70\u00f78=8, 6\",\n  \"97\u00f72=48, 1\",\n  \"35\u00f77=5, 0\",\n  \"94\u00f75=18, 4\",\n  \"79\u00f73=26, 1\",\n  \"75\u00f77=10, 5\",\n  \"13\u00f73=4, 1\",\n  \"48\u00f75=9, 3\",\n  \"62\u00f74=15, 2\",\n  \"72\u00f74=18, 0\",\n  \"77\u00f77=11, 0\",\n  \"94\u00f78=11, 6\",\n  \"14\u00f79=1, 5\",\n  \"54\u00f77=7, 5\",\n  \"90\u00f73=30, 0\",\n  \"78\u00f73=26, 0\",\n  \"92\u00f73=30, 2\",\n  \"28\u00f77=4, 0\",\n  \"11\u00f73=3, 2\",\n  \"93\u00f75=18, 3\",\n  \"46\u00f78=5, 6\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet next = 0;\nfor (let i = 0; i < paragraphs.items.length && next < replacements.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text && para.text.trim().length > 0) {\n    // Skip the leading date heading; only touch the division-answer cells.\n    if (/^\\s*\\d+\u00f7\\d+=\\d+,\\s*\\d+\\s*$/.test(para.text)) {\n      para.insertText(replacements[next], \"Replace\");\n      next++;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace the 25 division-problem answers in the (single) table, in\n# row-major document order, with their new values from the commit. Several\n# old answers are duplicated (e.g. \"14\u00f73=4, 2\" appears twice) but map to\n# different new values, so cells are addressed positionally (by walking the\n# table row-by-row, column-by-column and skipping the blank spacer rows)\n# rather than via text search/replace.\n\n$replacements = @(\n    \"34\u00f77=4, 6\",\n    \"86\u00f79=9, 5\",\n    \"90\u00f74=22, 2\",\n    \"85\u00f76=14, 1\",\n    \"70\u00f78=8, 6\",\n    \"97\u00f72=48, 1\",\n    \"35\u00f77=5, 0\",\n    \"94\u00f75=18, 4\",\n    \"79\u00f73=26, 1\",\n    \"75\u00f77=10, 5\",\n    \"13\u00f73=4, 1\",\n    \"48\u00f75=9, 3\",\n    \"62\u00f74=15, 2\",\n    \"72\u00f74=18, 0\",\n    \"77\u00f77=11, 0\",\n    \"94\u00f78=11, 6\",\n    \"14\u00f79=1, 5\",\n    \"54\u00f77=7, 5\",\n    \"90\u00f73=30, 0\",\n    \"78\u00f73=26, 0\",\n    \"92\u00f73=30, 2\",\n    \"28\u00f77=4, 0\",\n    \"11\u00f73=3, 2\",\n    \"93\u00f75=18, 3\",\n    \"46\u00f78=5, 6\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$next = 0\nfor ($r = 1; $r -le $t.Rows.Count; $r++) {\n    for ($c = 1; $c -le $t.Columns.Count; $c++) {\n        if ($next -ge $replacements.Length) { continue }\n        $cell = $t.Cell($r, $c)\n        $cellRange = $cell.Range\n        $plain = $cellRange.Text -replace \"[\\x07\\x0d]\", \"\"\n        if ($plain.Length -gt 0) {\n            $cellRange.Text = $replacements[$next]\n            $next = $next + 1\n        }\n    }\n}\n"}
